$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.652.51"
$ws.Range("E2").Value = "  -3.92%  "
$ws.Range("D3").Value = "3.367.44"
$ws.Range("E3").Value = "  -4.61%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'561.23"
$ws.Range("E5").Value = "  -3.63%  "
$ws.Range("D6").Value = "'183.26"
$ws.Range("E6").Value = "  -6.60%  "
$ws.Range("D7").Value = "'0.601"
$ws.Range("E7").Value = "  -1.54%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "3.359.27"
$ws.Range("E9").Value = "  -4.34%  "
$ws.Range("D10").Value = "'0.187"
$ws.Range("E10").Value = "  -7.82%  "
$ws.Range("D11").Value = "'0.595"
$ws.Range("E11").Value = "  -4.64%  "
$ws.Range("D12").Value = "'48.02"
$ws.Range("E12").Value = "  -7.05%  "
$ws.Range("D13").Value = "'0.0000270"
$ws.Range("E13").Value = "  -5.40%  "
$ws.Range("D14").Value = "'8.75"
$ws.Range("E14").Value = "  -5.47%  "
$ws.Range("D15").Value = "3.896.80"
$ws.Range("E15").Value = "  -4.63%  "
$ws.Range("D16").Value = "'607.41"
$ws.Range("E16").Value = "  -8.40%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "66.681.75"
$ws.Range("E17").Value = "  -3.94%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'18.16"
$ws.Range("E18").Value = "  -1.94%  "
$ws.Range("D19").Value = "3.361.20"
$ws.Range("E19").Value = "  -4.96%  "
$ws.Range("E20").Value = "  -3.14%  "
$ws.Range("D21").Value = "'11.53"
$ws.Range("E21").Value = "  -7.24%  "
$ws.Range("D22").Value = "'0.917"
$ws.Range("E22").Value = "  -5.09%  "
$ws.Range("D23").Value = "'16.89"
$ws.Range("E23").Value = "  -8.05%  "
$ws.Range("D24").Value = "'5.19"
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("D25").Value = "'100.34"
$ws.Range("E25").Value = "  -4.72%  "
$ws.Range("D26").Value = "'4.08"
$ws.Range("E26").Value = "  -6.60%  "
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").Value = "'2.72"
$ws.Range("E28").Value = "  -7.12%  "
$ws.Range("D29").Value = "'9.43"
$ws.Range("E29").Value = "  -7.33%  "
$ws.Range("D30").Value = "'8.82"
$ws.Range("E30").Value = "  -8.80%  "
$ws.Range("D31").Value = "'30.70"
$ws.Range("E31").Value = "  -7.95%  "
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").Value = "'3.86"
$ws.Range("E32").Value = "  -12.32%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'6.33"
$ws.Range("E33").Value = "  -7.49%  "
$ws.Range("D34").Value = "'11.16"
$ws.Range("E34").Value = "  -5.93%  "
$ws.Range("D35").Value = "'561.83"
$ws.Range("E35").Value = "  +11.51%  "
$ws.Range("D36").Value = "'0.106"
$ws.Range("E36").Value = "  -4.95%  "
$ws.Range("D37").Value = "3.849.28"
$ws.Range("E37").Value = "  +1.53%  "
$ws.Range("D38").Value = "'58.27"
$ws.Range("E38").Value = "  -5.89%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").Value = "'3.42"
$ws.Range("E40").Value = "  -6.84%  "
$ws.Range("D41").Value = "0.0₃0723"
$ws.Range("E41").Value = "  -10.95%  "
$ws.Range("D42").Value = "'3.47"
$ws.Range("E42").Value = "  +25.01%  "
$ws.Range("D43").Value = "'0.128"
$ws.Range("E43").Value = "  -4.43%  "
$ws.Range("D44").Value = "'2.67"
$ws.Range("E44").Value = "  -8.72%  "
$ws.Range("D45").Value = "'0.349"
$ws.Range("E45").Value = "  -6.20%  "
$ws.Range("D46").Value = "'32.36"
$ws.Range("E46").Value = "  -6.31%  "
$ws.Range("D47").Value = "'0.0417"
$ws.Range("E47").Value = "  -8.61%  "
$ws.Range("D48").Value = "'3.16"
$ws.Range("E48").Value = "  -6.74%  "
$ws.Range("D49").Value = "'2.67"
$ws.Range("E49").Value = "  -7.74%  "
$ws.Range("D50").Value = "'0.131"
$ws.Range("E50").Value = "  -3.97%  "
$ws.Range("D51").Value = "'0.998"
$ws.Range("E51").Value = "  -0.25%  "
